# LOM3043.docx edit: rotates several paragraph contents and restructures
# the "Avaliação" bullet paragraph (moves the "Método:" sentence out,
# relocates "Critério:" / "Norma de recuperação:" labels, and appends the
# bibliography text there), per the target diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simple paragraph-text swaps (content is rotated between sections).
#    These paragraphs have a single run each, so Range.Text assignment
#    is a safe, full replace.
# ---------------------------------------------------------------------

# "Objetivos" paragraph <- old "Programa resumido" text
$d.Paragraphs.Item(6).Range.Text = 'Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais.'

# "Docente(s)" bullet paragraph <- old "Objetivos" text
$d.Paragraphs.Item(8).Range.Text = 'Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia.'

# "Programa resumido" paragraph <- old "Programa" text
$d.Paragraphs.Item(10).Range.Text = '1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão.'

# "Programa" paragraph <- old "Método:" sentence (now a plain paragraph)
$d.Paragraphs.Item(12).Range.Text = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'

# "Bibliografia" paragraph <- old "Docente" text
$d.Paragraphs.Item(16).Range.Text = '5840622 - Miguel Justino Ribeiro Barboza'

# ---------------------------------------------------------------------
# 2) Rebuild the "Avaliação" paragraph (item 14). It has several runs
#    (bold labels + plain text), so Range.Text assignment would only
#    touch the first run; instead we delete its contents (keeping the
#    paragraph mark) and retype it, then re-bold the three labels.
#
#    New layout:
#      Método:  | A média do semestre será computada com base na relação: <br><br> M=(P1+2P2)/3 <br>
#      Critério: | A recuperação será composta por uma única prova (RC) ... <br><br>
#                  A média final, para os alunos em recuperação, ... <br><br> MF=(M+RC)/2 <br>
#      Norma de recuperação: | <bibliography text>
# ---------------------------------------------------------------------

$avalParagraph = $d.Paragraphs.Item(14)
$body = $avalParagraph.Range.Duplicate
$body.MoveEnd(1, -1)          # exclude the trailing paragraph mark
$body.Delete()

$LF = [char]11                # maps to <w:br/> when typed into a run

$full = (
    'Método: ' +
    'A média do semestre será computada com base na relação:' + $LF + $LF +
    'M=(P1+2P2)/3' + $LF +
    'Critério: ' +
    'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.' + $LF + $LF +
    'A média final, para os alunos em recuperação, será computada com base na relação abaixo:' + $LF + $LF +
    'MF=(M+RC)/2' + $LF +
    'Norma de recuperação: ' +
    '1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006.'
)

$insPoint = $d.Paragraphs.Item(14).Range.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertAfter($full)

# Re-apply bold formatting to the three labels within the rebuilt paragraph.
foreach ($label in @('Método: ', 'Critério: ', 'Norma de recuperação: ')) {
    $scope = $d.Paragraphs.Item(14).Range.Duplicate
    $scope.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    $scope.Font.Bold = 1
}
